$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 20.32821866666667
$ws.Range("H2").Value = 60.984656
$ws.Range("I2").Value = 0.004181898474048532
$ws.Range("J2").Value = 0.004181898474048532
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.82741333333333
$ws.Range("N2").Value = 95.48223999999999
$ws.Range("O2").Value = 0.114390792932228
$ws.Range("P2").Value = 0.114390792932228
$ws.Range("Q2").Value = 646.9946178343822
$ws.Range("R2").Value = 5822.951560509439
$ws.Range("S2").Value = 0.0004783706824084859
$ws.Range("T2").Value = 0.000478370682408486
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 20.32821866666667
$ws.Range("H3").Value = 60.984656
$ws.Range("I3").Value = 0.004181898474048532
$ws.Range("J3").Value = 0.004181898474048532
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 85.46317833333335
$ws.Range("N3").Value = 256.389535
$ws.Range("O3").Value = 0.307162904935779
$ws.Range("P3").Value = 0.307162904935779
$ws.Range("Q3").Value = 1737.314177108329
$ws.Range("R3").Value = 15635.82759397496
$ws.Range("S3").Value = 0.001284524083435249
$ws.Range("T3").Value = 0.001284524083435249
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 20.32821866666667
$ws.Range("H4").Value = 60.984656
$ws.Range("I4").Value = 0.004181898474048532
$ws.Range("J4").Value = 0.004181898474048532
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 122.2478306666667
$ws.Range("N4").Value = 366.743492
$ws.Range("O4").Value = 0.4393704929064738
$ws.Range("P4").Value = 0.4393704929064738
$ws.Range("Q4").Value = 2485.080633317639
$ws.Range("R4").Value = 22365.72569985875
$ws.Range("S4").Value = 0.001837402793827534
$ws.Range("T4").Value = 0.001837402793827534
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 20.32821866666667
$ws.Range("H5").Value = 60.984656
$ws.Range("I5").Value = 0.004181898474048532
$ws.Range("J5").Value = 0.004181898474048532
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 38.69562533333333
$ws.Range("N5").Value = 116.086876
$ws.Range("O5").Value = 0.1390758092255191
$ws.Range("P5").Value = 0.1390758092255191
$ws.Range("Q5").Value = 786.6131332194062
$ws.Range("R5").Value = 7079.518198974655
$ws.Range("S5").Value = 0.0005816009143772632
$ws.Range("T5").Value = 0.0005816009143772632
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4809.896321333334
$ws.Range("H6").Value = 14429.688964
$ws.Range("I6").Value = 0.9894865072215304
$ws.Range("J6").Value = 0.9894865072215304
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.82741333333333
$ws.Range("N6").Value = 95.48223999999999
$ws.Range("O6").Value = 0.114390792932228
$ws.Range("P6").Value = 0.114390792932228
$ws.Range("Q6").Value = 153086.5583095555
$ws.Range("R6").Value = 1377779.024785999
$ws.Range("S6").Value = 0.1131881461568116
$ws.Range("T6").Value = 0.1131881461568117
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4809.896321333334
$ws.Range("H7").Value = 14429.688964
$ws.Range("I7").Value = 0.9894865072215304
$ws.Range("J7").Value = 0.9894865072215304
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 85.46317833333335
$ws.Range("N7").Value = 256.389535
$ws.Range("O7").Value = 0.307162904935779
$ws.Range("P7").Value = 0.307162904935779
$ws.Range("Q7").Value = 411069.0270749547
$ws.Range("R7").Value = 3699621.243674592
$ws.Range("S7").Value = 0.303933549952923
$ws.Range("T7").Value = 0.303933549952923
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4809.896321333334
$ws.Range("H8").Value = 14429.688964
$ws.Range("I8").Value = 0.9894865072215304
$ws.Range("J8").Value = 0.9894865072215304
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 122.2478306666667
$ws.Range("N8").Value = 366.743492
$ws.Range("O8").Value = 0.4393704929064738
$ws.Range("P8").Value = 0.4393704929064738
$ws.Range("Q8").Value = 587999.3910145803
$ws.Range("R8").Value = 5291994.519131223
$ws.Range("S8").Value = 0.4347511744022289
$ws.Range("T8").Value = 0.434751174402229
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4809.896321333334
$ws.Range("H9").Value = 14429.688964
$ws.Range("I9").Value = 0.9894865072215304
$ws.Range("J9").Value = 0.9894865072215304
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 38.69562533333333
$ws.Range("N9").Value = 116.086876
$ws.Range("O9").Value = 0.1390758092255191
$ws.Range("P9").Value = 0.1390758092255191
$ws.Range("Q9").Value = 186121.945942493
$ws.Range("R9").Value = 1675097.513482436
$ws.Range("S9").Value = 0.1376136367095668
$ws.Range("T9").Value = 0.1376136367095668
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.69506
$ws.Range("H10").Value = 8.085180000000001
$ws.Range("I10").Value = 0.000554424737665286
$ws.Range("J10").Value = 0.000554424737665286
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.82741333333333
$ws.Range("N10").Value = 95.48223999999999
$ws.Range("O10").Value = 0.114390792932228
$ws.Range("P10").Value = 0.114390792932228
$ws.Range("Q10").Value = 85.77678857813333
$ws.Range("R10").Value = 771.9910972032001
$ws.Range("S10").Value = 0.00006342108536277458
$ws.Range("T10").Value = 0.00006342108536277459
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.69506
$ws.Range("H11").Value = 8.085180000000001
$ws.Range("I11").Value = 0.000554424737665286
$ws.Range("J11").Value = 0.000554424737665286
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 85.46317833333335
$ws.Range("N11").Value = 256.389535
$ws.Range("O11").Value = 0.307162904935779
$ws.Range("P11").Value = 0.307162904935779
$ws.Range("Q11").Value = 230.3283933990334
$ws.Range("R11").Value = 2072.955540591301
$ws.Range("S11").Value = 0.0001702987129895265
$ws.Range("T11").Value = 0.0001702987129895265
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.69506
$ws.Range("H12").Value = 8.085180000000001
$ws.Range("I12").Value = 0.000554424737665286
$ws.Range("J12").Value = 0.000554424737665286
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 122.2478306666667
$ws.Range("N12").Value = 366.743492
$ws.Range("O12").Value = 0.4393704929064738
$ws.Range("P12").Value = 0.4393704929064738
$ws.Range("Q12").Value = 329.4652385165067
$ws.Range("R12").Value = 2965.18714664856
$ws.Range("S12").Value = 0.0002435978702675391
$ws.Range("T12").Value = 0.0002435978702675392
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.69506
$ws.Range("H13").Value = 8.085180000000001
$ws.Range("I13").Value = 0.000554424737665286
$ws.Range("J13").Value = 0.000554424737665286
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 38.69562533333333
$ws.Range("N13").Value = 116.086876
$ws.Range("O13").Value = 0.1390758092255191
$ws.Range("P13").Value = 0.1390758092255191
$ws.Range("Q13").Value = 104.2870320108533
$ws.Range("R13").Value = 938.58328809768
$ws.Range("S13").Value = 0.00007710706904544581
$ws.Range("T13").Value = 0.00007710706904544581
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 28.08283533333333
$ws.Range("H14").Value = 84.24850599999999
$ws.Range("I14").Value = 0.005777169566755752
$ws.Range("J14").Value = 0.005777169566755752
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 31.82741333333333
$ws.Range("N14").Value = 95.48223999999999
$ws.Range("O14").Value = 0.114390792932228
$ws.Range("P14").Value = 0.114390792932228
$ws.Range("Q14").Value = 893.8040077259376
$ws.Range("R14").Value = 8044.236069533438
$ws.Range("S14").Value = 0.0006608550076451266
$ws.Range("T14").Value = 0.0006608550076451268
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 28.08283533333333
$ws.Range("H15").Value = 84.24850599999999
$ws.Range("I15").Value = 0.005777169566755752
$ws.Range("J15").Value = 0.005777169566755752
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 85.46317833333335
$ws.Range("N15").Value = 256.389535
$ws.Range("O15").Value = 0.307162904935779
$ws.Range("P15").Value = 0.307162904935779
$ws.Range("Q15").Value = 2400.048364198301
$ws.Range("R15").Value = 21600.43527778471
$ws.Range("S15").Value = 0.001774532186431273
$ws.Range("T15").Value = 0.001774532186431273
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 28.08283533333333
$ws.Range("H16").Value = 84.24850599999999
$ws.Range("I16").Value = 0.005777169566755752
$ws.Range("J16").Value = 0.005777169566755752
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 122.2478306666667
$ws.Range("N16").Value = 366.743492
$ws.Range("O16").Value = 0.4393704929064738
$ws.Range("P16").Value = 0.4393704929064738
$ws.Range("Q16").Value = 3433.065698469217
$ws.Range("R16").Value = 30897.59128622295
$ws.Range("S16").Value = 0.002538317840149754
$ws.Range("T16").Value = 0.002538317840149755
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 28.08283533333333
$ws.Range("H17").Value = 84.24850599999999
$ws.Range("I17").Value = 0.005777169566755752
$ws.Range("J17").Value = 0.005777169566755752
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 38.69562533333333
$ws.Range("N17").Value = 116.086876
$ws.Range("O17").Value = 0.1390758092255191
$ws.Range("P17").Value = 0.1390758092255191
$ws.Range("Q17").Value = 1086.682874356362
$ws.Range("R17").Value = 9780.145869207254
$ws.Range("S17").Value = 0.0008034645325295979
$ws.Range("T17").Value = 0.0008034645325295979
